$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.185770869255066
$ws.Range("B1").Value = 2.340942144393921
$ws.Range("C1").Value = 5.004730224609375
$ws.Range("D1").Value = 2.416975259780884
$ws.Range("E1").Value = 1.219061255455017
